$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the model-parameter columns (BQ:BV). Editing the header cells keeps
# the backing Excel Table ("Table3") column names in sync automatically.
$ws.Range("BQ1").Value = "al_win"
$ws.Range("BR1").Value = "al_loss"
$ws.Range("BS1").Value = "bet_win"
$ws.Range("BT1").Value = "bet_loss"
$ws.Range("BU1").Value = "sc_bet_win"
$ws.Range("BV1").Value = "sc_bet_loss"

# Refresh the fitted-model parameter values (best model refit) for BQ:BV, rows 2-29.
$modelParams = @{
    2 = @(0.033749659017537598, 0.34601990922606202, 1.74152919562174, 0.47967486629865602, -0.48246735807442798, -0.346096577956584)
    3 = @(0.22711067812746599, 0.57114357630567403, 1.10546618297845, 0.54185088284372296, -0.149792832644881, 0.28552077400074899)
    4 = @(0.168468865392258, 0.29856914159399001, 1.46099372204197, 1.0337159556990101, 0.13155603791655601, -0.046111333073499003)
    5 = @(0.143532833929191, 0.30642130850684601, 1.9133346438900201, 1.31073622608664, 0.141362091414032, 0.50482723298507604)
    6 = @(0.090379433712604498, 0.70276377885577701, 0.91939584489317205, 1.02683493665536, 0.25227592851584801, -0.28584661775709702)
    7 = @(0.108149850154973, 0.28373926618510298, 1.84186107435012, 1.6404560698132, -0.15455867893363101, 0.23331283688676599)
    8 = @(0.156593073875155, 0.39302413997605301, 1.9452254807627001, 1.62732142365252, 0.31154256063091901, -0.40179045698015903)
    9 = @(0.29651199330697797, 0.64135026446160803, 0.82787544850094497, 0.13453166119818599, 0.72999667581427197, 1.40144688210716)
    10 = @(0.107565932900949, 0.13500411064186699, 2.94487344515893, 1.6123714771535, 0.17348764213457499, 1.2580357289008799)
    11 = @(0.197304510781806, 0.28584479782841699, 0.58288409722016599, 0.52607110064336404, 0.70853159514417996, 0.49921830533690498)
    12 = @(0.443162596624875, 0.51995011929615098, 1.13982551054293, 0.99202346820181997, -0.15846268473448699, 0.16787458272271299)
    13 = @(0.215528105305304, 0.30652054875453399, 2.1839780414740102, 1.94047625893135, -0.060761536405578402, -0.81953891198047701)
    14 = @(0.45457115321685698, 0.50571942681840698, 1.35335469809391, 0.75650556699043503, 0.17011508376602, 0.37579360453347299)
    15 = @(0.373434423478611, 0.460098826390475, 1.5306011420310599, 0.56547843668977704, 0.078578762001784805, -0.84655485837810596)
    16 = @(0.19604022275664701, 0.35035087748668298, 2.3922539209142499, 1.88571675302126, -0.27389114201592901, -0.35290988002387302)
    17 = @(0.10209777107449899, 0.38112066674259498, 1.68348988229774, 0.93049963251203904, 0.23762322213363499, 0.0498494571118904)
    18 = @(0.188924251774887, 0.66270029189331303, 1.6210445900649, 1.2784920597824001, 0.14197864775852201, -0.13141936209968499)
    19 = @(0.0147408103222445, 0.100391849746247, 2.1119107178987999, 0.87684533040774404, -0.13404313871812601, -0.38127798850820399)
    20 = @(0.241294417481065, 0.43362811489248598, 2.12083714768929, 1.2602974438961401, -0.480528276238676, -0.19015492769551801)
    21 = @(0.33289019028123701, 0.442045457138729, 1.1537117661419101, 0.54437650501970403, 0.13735357715461599, -0.30504896353712202)
    22 = @(0.20564889429056299, 0.40186403845827301, 1.79293778213238, 1.43303948217488, 0.079673160204535201, -0.97505363676711299)
    23 = @(0.214718034975173, 0.44823018226292399, 1.7953024570010201, 1.06712222423887, -0.53540111260789902, 0.31007950480854701)
    24 = @(0.042428203848660098, 0.063980076271133099, 1.83228156649898, 1.7934711077149701, 0.037947630266227002, -0.77441235839064204)
    25 = @(0.20403438305475999, 0.23610185335377301, 2.0239735182112102, 2.3314785751660398, -0.20121067857813299, -0.68347323761426704)
    26 = @(0.140533382413114, 0.221742667078277, 1.84442343638814, 2.3879442920877598, -0.062232260737934003, -0.32530782700728)
    27 = @(0.23282562272131299, 0.65056803878946301, 1.7782001530629601, 0.51214732272857799, -0.17559169365990401, -0.72686406206242404)
    28 = @(0.16351552882031301, 0.214621949934555, 2.8685449948128801, 1.62059239386302, 0.10024782188338099, 0.35687285947530001)
    29 = @(0.070548161359247699, 0.22567393990067999, 2.1874361525966601, 1.8522920520388799, 0.45147420617223999, 0.13403400562437701)
}

foreach ($row in $modelParams.Keys) {
    $vals = $modelParams[$row]
    $ws.Cells.Item($row, 69).Value = $vals[0]
    $ws.Cells.Item($row, 70).Value = $vals[1]
    $ws.Cells.Item($row, 71).Value = $vals[2]
    $ws.Cells.Item($row, 72).Value = $vals[3]
    $ws.Cells.Item($row, 73).Value = $vals[4]
    $ws.Cells.Item($row, 74).Value = $vals[5]
}

# Move the active selection to reflect where the user left off reviewing the new columns.
$null = $ws.Range("CA10").Select()
